# Update biodiversity indicators for age and species (solid_wood / sum_product / ecosystem / system)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ B = 1146999.278263127; D = 1146999.278263127; E = 49458556.73467615; F = 50605556.01293928 }
    3  = @{ B = 1525165.880659537; D = 1525165.880659537; E = 48485031.07559408; F = 50010196.95625361 }
    4  = @{ B = 1819700.853357395; D = 1819700.853357395; E = 47982073.2953876;  F = 49801774.14874499 }
    5  = @{ B = 1901179.583824972; D = 1901179.583824972; E = 47931184.59034887; F = 49832364.17417384 }
    6  = @{ B = 1943462.411646292; D = 1943462.411646292; E = 47894441.1993427;  F = 49837903.61098899 }
    7  = @{ B = 1854427.755306635; D = 1854427.755306635; E = 47785454.00699076; F = 49639881.76229739 }
    8  = @{ B = 1944765.397619904; D = 1944765.397619904; E = 48081313.03371169; F = 50026078.43133159 }
    9  = @{ B = 1618751.847318216; D = 1618751.847318216; E = 48395835.40924136; F = 50014587.25655958 }
    10 = @{ B = 1281538.70612917;  D = 1281538.70612917;  E = 49234707.21529362; F = 50516245.92142279 }
    11 = @{ B = 1647603.39473688;  D = 1647603.39473688;  E = 50711989.88827336; F = 52359593.28301024 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("B$row").Value = $rowVals.B
    $ws.Range("D$row").Value = $rowVals.D
    $ws.Range("E$row").Value = $rowVals.E
    $ws.Range("F$row").Value = $rowVals.F
}
